$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "Question No"
$ws.Range("B1").Value = "GFG/LC"
$ws.Range("C1").Value = "Question"

# ---- Data rows ----
$ws.Range("A2").Value = 94
$ws.Range("B2").Value = "LC"
$ws.Range("C2").Value = "Binary Tree Inorder Traversal(Inorder Tree Traversal -Recursive)"

$ws.Range("A3").Value = 145
$ws.Range("B3").Value = "LC"
$ws.Range("C3").Value = "Binary Tree Postorder Traversal-Recursive"

$ws.Range("A4").Value = 144
$ws.Range("B4").Value = "LC"
$ws.Range("C4").Value = "Binary Tree Preorder Traversal-Recursive"

$ws.Range("A5").Value = 701
$ws.Range("B5").Value = "LC"
$ws.Range("C5").Value = "Insert into a Binary Search Tree"

$ws.Range("A6").Value = 700
$ws.Range("B6").Value = "LC"
$ws.Range("C6").Value = "Search in a Binary Search Tree"

$ws.Range("A7").Value = 216
$ws.Range("B7").Value = "GFG"
$ws.Range("C7").Value = "Minimum element in BST (Find min and max value in a BST)"

# ---- Column C width (~55 chars once round-tripped through Excel) ----
$ws.Columns.Item(3).ColumnWidth = 54.16666666666667

# ---- Row heights for the two wrapped header/first-data rows ----
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 30

# ---- Build body alignment (left/top + wrap) on a single template cell, then
#      fan it out with a format-only paste so we don't churn the style table ----
$bodyTemplate = $ws.Range("A2")
$bodyTemplate.HorizontalAlignment = -4131   # xlLeft
$bodyTemplate.VerticalAlignment = -4160     # xlTop
$bodyTemplate.WrapText = $true

$bodyTemplate.Copy()
$ws.Range("A1:C7").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Header formatting: bold + left/top + wrap, built on a single template cell ----
$headerTemplate = $ws.Range("A1")
$headerTemplate.Font.Bold = $true
$headerTemplate.HorizontalAlignment = -4131
$headerTemplate.VerticalAlignment = -4160
$headerTemplate.WrapText = $true

$headerTemplate.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- "Question" header cell (C1): bold + left/top, no wrap ----
$questionHeader = $ws.Range("C1")
$questionHeader.Font.Bold = $true
$questionHeader.HorizontalAlignment = -4131
$questionHeader.VerticalAlignment = -4160
$questionHeader.WrapText = $false

# ---- Page setup ----
$ws.PageSetup.Orientation = 1   # xlPortrait

# ---- Selection ----
$ws.Range("C7").Select()

$wb.Save()
